$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the refreshed crypto feed.
# D-column values are numeric-looking text (e.g. "309.96"); Excel would silently
# coerce a bare .Value assignment into a float (and mutate the cell style via an
# explicit NumberFormat). Force text with NumberFormat="@" then restore the
# original (default/no) style so only the value itself changes, matching the
# upstream inline-string cells exactly.

$origStyle = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.956.12'
$ws.Range('D2').Style = $origStyle
$ws.Range('E2').Value = '  +0.33%  '
$origStyle = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.818.57'
$ws.Range('D3').Style = $origStyle
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  +0.15%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '309.96'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('E6').Value = '  +0.13%  '
$origStyle = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4653'
$ws.Range('D7').Style = $origStyle
$ws.Range('E7').Value = '  +0.23%  '
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3662'
$ws.Range('D8').Style = $origStyle
$ws.Range('E8').Value = '  -1.10%  '
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07360'
$ws.Range('D9').Style = $origStyle
$ws.Range('E9').Value = '  +0.07%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.8720'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  -0.65%  '
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.25'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  -1.11%  '
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.824.21'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  +3.31%  '
$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.403'
$ws.Range('D13').Style = $origStyle
$ws.Range('E13').Value = '  +0.78%  '
$ws.Range('E14').Value = '  +0.92%  '
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.510'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  -0.15%  '
$origStyle = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '91.48'
$ws.Range('D16').Style = $origStyle
$ws.Range('E16').Value = '  -0.43%  '
$ws.Range('E17').Value = '  +0.19%  '
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008703'
$ws.Range('D18').Style = $origStyle
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('E20').Value = '  -0.66%  '
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '26.972.02'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  +0.39%  '
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.295'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  -0.38%  '
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.58'
$ws.Range('D23').Style = $origStyle
$ws.Range('E23').Value = '  +0.04%  '
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.048.15'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  +1.89%  '
$ws.Range('E25').Value = '  -0.08%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '150.69'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  -0.64%  '
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.36'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  -0.22%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.137'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  -0.95%  '
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.254'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  -1.42%  '
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '116.52'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  +0.32%  '
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08897'
$ws.Range('D31').Style = $origStyle
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.7581'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  +0.61%  '
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.164'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  +0.50%  '
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.503'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  +0.83%  '
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.905'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  -0.79%  '
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('E37').Value = '  -1.34%  '
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.05290'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  +0.74%  '
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01947'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  -1.01%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.972'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  +1.51%  '
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '7.182'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  +0.09%  '
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.5284'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  -0.81%  '
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.352'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  -2.76%  '
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.1659'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  -0.36%  '
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.442'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  -0.70%  '
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4872'
$ws.Range('D46').Style = $origStyle
$ws.Range('E46').Value = '  -2.22%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.46'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  +1.01%  '
$ws.Range('E48').Value = '  +0.15%  '
$ws.Range('E49').Value = '  -0.29%  '
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '103.36'
$ws.Range('D50').Style = $origStyle
$ws.Range('E50').Value = '  -0.34%  '
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06296'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  +0.01%  '
